# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with refreshed values. Price-column cells are assigned with a
# leading apostrophe to force Excel to store them as text (otherwise
# decimal-looking strings like "405.78" would be coerced to numbers and
# lose trailing zeros / gain floating-point noise); the Style is then
# reset to "Normal" so the cell doesn't retain the quote-prefix style
# Excel applies automatically, keeping formatting identical to before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'61.221.42"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -0.48%  '
$ws.Cells.Item(3, 4).Value = "'3.376.73"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -0.88%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).Value = "'405.78"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -1.67%  '
$ws.Cells.Item(6, 4).Value = "'135.02"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +8.97%  '
$ws.Cells.Item(7, 5).Value = '  +1.24%  '
$ws.Cells.Item(8, 5).Value = '  +0.03%  '
$ws.Cells.Item(9, 4).Value = "'0.675"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +5.24%  '
$ws.Cells.Item(10, 5).Value = '  -2.61%  '
$ws.Cells.Item(11, 4).Value = "'42.86"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +4.03%  '
$ws.Cells.Item(12, 5).Value = '  -0.80%  '
$ws.Cells.Item(13, 4).Value = "'3.889.04"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -1.50%  '
$ws.Cells.Item(14, 4).Value = "'8.37"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -0.60%  '
$ws.Cells.Item(15, 4).Value = "'19.61"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(16, 4).Value = "'3.383.30"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.72%  '
$ws.Cells.Item(17, 4).Value = "'61.132.14"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.51%  '
$ws.Cells.Item(18, 5).Value = '  -0.55%  '
$ws.Cells.Item(19, 4).Value = "'11.04"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +1.81%  '
$ws.Cells.Item(20, 5).Value = '  +3.50%  '
$ws.Cells.Item(21, 5).Value = '  -2.91%  '
$ws.Cells.Item(22, 4).Value = "'83.60"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +9.13%  '
$ws.Cells.Item(23, 4).Value = "'314.16"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +5.17%  '
$ws.Cells.Item(24, 4).Value = "'12.75"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -0.52%  '
$ws.Cells.Item(25, 5).Value = '  -0.46%  '
$ws.Cells.Item(26, 4).Value = "'4.77"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +11.59%  '
$ws.Cells.Item(27, 4).Value = "'8.37"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +10.04%  '
$ws.Cells.Item(28, 5).Value = '  -3.54%  '
$ws.Cells.Item(29, 4).Value = "'7.47"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -7.11%  '
$ws.Cells.Item(30, 5).Value = '  +0.50%  '
$ws.Cells.Item(31, 4).Value = "'0.117"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.00%  '
$ws.Cells.Item(32, 4).Value = "'1.00"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +0.01%  '
$ws.Cells.Item(33, 4).Value = "'11.32"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -0.73%  '
$ws.Cells.Item(34, 4).Value = "'41.36"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -2.86%  '
$ws.Cells.Item(35, 5).Value = '  -2.29%  '
$ws.Cells.Item(36, 5).Value = '  +0.62%  '
$ws.Cells.Item(37, 4).Value = "'52.23"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -0.37%  '
$ws.Cells.Item(38, 4).Value = "'0.996"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -0.13%  '
$ws.Cells.Item(39, 4).Value = "'3.42"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -2.61%  '
$ws.Cells.Item(40, 5).Value = '  -2.88%  '
$ws.Cells.Item(41, 4).Value = "'137.62"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +2.38%  '
$ws.Cells.Item(42, 5).Value = '  +1.08%  '
$ws.Cells.Item(43, 4).Value = "'0.124"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.62%  '
$ws.Cells.Item(44, 5).Value = '  +3.24%  '
$ws.Cells.Item(45, 5).Value = '  +3.05%  '
$ws.Cells.Item(46, 4).Value = "'16.65"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -3.74%  '
$ws.Cells.Item(47, 4).Value = "'2.23"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +2.04%  '
$ws.Cells.Item(48, 4).Value = "'21.34"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -2.45%  '
$ws.Cells.Item(49, 4).Value = "'2.125.71"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -3.34%  '
$ws.Cells.Item(50, 5).Value = '  -5.37%  '
$ws.Cells.Item(51, 5).Value = '  -0.19%  '
